# [FEATURE] Add filter option list
# Adds three new Arbeitsmatrix rows (new issues) plus two blank spacer rows
# right before the "Stunden insgesamt" summary row, which shifts from row
# 105 down to row 110.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Arbeitsmatrix")
$ws.Activate()

# ---------------------------------------------------------------------
# 1) Make room: insert 5 new rows right before the old summary row (105).
#    This pushes the existing summary row (with its SUM(I:I)/SUM(H:H)
#    formulas) down to row 110 automatically, untouched.
# ---------------------------------------------------------------------
$ws.Rows("105:109").Insert()

# ---------------------------------------------------------------------
# 2) Copy the formatting of two existing "highlighted" data rows (style
#    index pattern used for this issue block) onto the two new data rows
#    105 and 106.
# ---------------------------------------------------------------------
$ws.Range("A103:K103").Copy()
$ws.Range("A105:K106").PasteSpecial(-4122)  # xlPasteFormats

# Row 107 uses the plain / unstyled look (like row 97) for A, B, C, E.
# D107, F107, G107, I107, J107, K107 keep the same look as the row above
# (pasted first so the later, narrower E107 paste is not clobbered).
$ws.Range("D103").Copy()
$ws.Range("D107").PasteSpecial(-4122)
$ws.Range("F103:K103").Copy()
$ws.Range("F107:K107").PasteSpecial(-4122)

$ws.Range("A97:C97").Copy()
$ws.Range("A107:C107").PasteSpecial(-4122)
$ws.Range("E97").Copy()
$ws.Range("E107").PasteSpecial(-4122)

$excel.CutCopyMode = 0

# The template rows (103) carry a formatted-but-empty "H" (Stunden Seminar)
# cell; the new rows never had that column touched, so drop it completely.
$ws.Range("H105:H107").Clear()

# ---------------------------------------------------------------------
# 3) Blank spacer rows 108 & 109 only need formatting on D/F/G, matching
#    the look of the other spacer rows in the sheet (e.g. row 95).
# ---------------------------------------------------------------------
foreach ($r in 108..109) {
    $ws.Cells.Item(95, 4).Copy()
    $ws.Cells.Item($r, 4).PasteSpecial(-4122)
    $ws.Cells.Item(95, 6).Copy()
    $ws.Cells.Item($r, 6).PasteSpecial(-4122)
    $ws.Cells.Item(95, 7).Copy()
    $ws.Cells.Item($r, 7).PasteSpecial(-4122)
}
$excel.CutCopyMode = 0

# ---------------------------------------------------------------------
# 4) Fill in the new issue data. Enter the new shared-string text in the
#    same order the author originally typed it in, so the shared string
#    table keeps the same ordering (127=Alle Filter.., 128=Icons bauen,
#    129=Rezept Teaser).
# ---------------------------------------------------------------------
$ws.Cells.Item(107, 5).Value = "Alle Filter Optionen Sammeln"
$ws.Cells.Item(105, 5).Value = "MockUps Icons bauen"
$ws.Cells.Item(106, 5).Value = "MockUps Rezept Teaser"

# Row 105 - Interface Design / MockUps / [FEATURE]
$ws.Cells.Item(105, 1).Value = 22
$ws.Cells.Item(105, 2).Value = "Interface Design"
$ws.Cells.Item(105, 3).Value = "MockUps"
$ws.Cells.Item(105, 4).Value = "[FEATURE]"
$ws.Cells.Item(105, 6).Value = 44375
$ws.Cells.Item(105, 7).Value = 44359
$ws.Cells.Item(105, 9).Formula = "=ROUNDUP(((SUM(K105-J105)*24*60/60)/0.25),0)*0.25"
$ws.Cells.Item(105, 10).Value = 0.625
$ws.Cells.Item(105, 11).Value = 0.79166666666666663

# Row 106 - Interface Design / MockUps / [FEATURE]
$ws.Cells.Item(106, 1).Value = 22
$ws.Cells.Item(106, 2).Value = "Interface Design"
$ws.Cells.Item(106, 3).Value = "MockUps"
$ws.Cells.Item(106, 4).Value = "[FEATURE]"
$ws.Cells.Item(106, 6).Value = 44376
$ws.Cells.Item(106, 7).Value = 44359
$ws.Cells.Item(106, 9).Formula = "=ROUNDUP(((SUM(K106-J106)*24*60/60)/0.25),0)*0.25"
$ws.Cells.Item(106, 10).Value = 0.375
$ws.Cells.Item(106, 11).Value = 0.625

# Row 107 - Konzeptuelles Design / Content Map / [FEATURE]
$ws.Cells.Item(107, 1).Value = 18
$ws.Cells.Item(107, 2).Value = "Konzeptuelles Design"
$ws.Cells.Item(107, 3).Value = "Content Map"
$ws.Cells.Item(107, 4).Value = "[FEATURE]"
$ws.Cells.Item(107, 6).Value = 44376
$ws.Cells.Item(107, 7).Value = 44359
$ws.Cells.Item(107, 9).Formula = "=ROUNDUP(((SUM(K107-J107)*24*60/60)/0.25),0)*0.25"
$ws.Cells.Item(107, 10).Value = 0.66666666666666663
$ws.Cells.Item(107, 11).Value = 0.70833333333333337

# ---------------------------------------------------------------------
# 5) Extend the "Prefix" dropdown validation (the $N$3:$N$6 list) to the
#    newly added rows (it previously ended at D104).
# ---------------------------------------------------------------------
$validationRange = $ws.Range("D105:D109")
$validationRange.Validation.Add(3, 1, 1, "=`$N`$3:`$N`$6")
$validationRange.Validation.ErrorTitle = "Prefix nicht unterstützt"
$validationRange.Validation.ErrorMessage = "Es konnte kein korrekter Prefix ausgegeben werden`n"
$validationRange.Validation.InputTitle = "Prefix"
$validationRange.Validation.InputMessage = "Wählen Sie einen Prefix aus"
$validationRange.Validation.IgnoreBlank = $true
$validationRange.Validation.InCellDropdown = $true
$validationRange.Validation.ShowInput = $true
$validationRange.Validation.ShowError = $true

# ---------------------------------------------------------------------
# 6) Update the view: scroll/selection now centers on the new rows.
# ---------------------------------------------------------------------
$ws.Range("H110").Select() | Out-Null
